$wb = $excel.ActiveWorkbook

# --- Indicator sheet: fill in field "Occurrence" types for B3:B23 ---
$ws = $wb.Worksheets.Item("Indicator")

$ws.Range("B3").Value = "required"
$ws.Range("B4").Value = "optional"
$ws.Range("B5").Value = "required"
$ws.Range("B6").Value = "required"
$ws.Range("B7").Value = "optional"
$ws.Range("B8").Value = "optional"
$ws.Range("B9").Value = "optional"
$ws.Range("B10").Value = "required"
$ws.Range("B11").Value = "optional"
$ws.Range("B12").Value = "optional"
$ws.Range("B13").Value = "optional"
$ws.Range("B14").Value = "prohibited"
$ws.Range("B15").Value = "optional"
$ws.Range("B16").Value = "optional"
$ws.Range("B17").Value = "prohibited"
$ws.Range("B18").Value = "optional"
$ws.Range("B19").Value = "optional"
$ws.Range("B20").Value = "optional"
$ws.Range("B21").Value = "prohibited"
$ws.Range("B22").Value = "prohibited"
$ws.Range("B23").Value = "optional"

# --- Update selection on the TTP sheet without making it the active tab ---
$ttp = $wb.Worksheets.Item("TTP")
$ttp.Range("D24").Select() | Out-Null

# --- Make Indicator the active sheet/tab and update its selection ---
$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
